$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for column G (K) for rows 2-25, replacing the old Strike# values
$newValues = @{
    2  = 3
    3  = 2
    4  = 0
    5  = 4
    6  = 1
    7  = 0
    8  = 1
    9  = 0
    10 = 1
    11 = 2
    12 = 0
    13 = 2
    14 = 2
    15 = 0
    16 = 2
    17 = 2
    18 = 3
    19 = 0
    20 = 1
    21 = 2
    22 = 3
    23 = 3
    24 = 1
    25 = 1
}

foreach ($row in $newValues.Keys) {
    $ws.Range("G$row").Value = $newValues[$row]
}
